$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 75; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 46060) {
        $cell.Value = 46061
    }
}
